$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = "25.48."
$ws.Range("G2").Value = "25.48."
$ws.Range("G3").Value = "21.09."
$ws.Range("G4").Value = "17.52."
